$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUD")

# Insert a new column before AE (col 31): shifts AE..AI -> AF..AJ
$ws.Columns("AE:AE").Insert()

# New column header (row 13, merged 12:14 area) and value (row 15)
$ws.Range("AE13").Value = "Autre allocation"
$ws.Range("AE15").Value = 20000

# Fix up the named ranges that the column insert should have grown
$wb.Names.Item("SUD!_FilterDatabase").RefersTo = "=SUD!`$A`$14:`$AK`$16"
$wb.Names.Item("SUD!Print_Area").RefersTo = "=SUD!`$A`$1:`$AG`$15"

# Restore selection / view state on the SUD sheet
$ws.Activate()
$ws.Range("AE16").Select()
